# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp string
# - Update COVID-19 case counters for several countries
# - Filipinas overtakes Indonesia in case totals (rows 26/27 swap identities)
# - Estonia overtakes Ruanda in case totals (rows 128/129 swap identities)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 10:39"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 4918770
$ws.Range("C4").Value = 350
$ws.Range("D4").Value = 2482899
$ws.Range("E4").Value = 2275553
$ws.Range("G4").Value = 28
$ws.Range("H4").Value = 160318

# Rusia (row 7)
$ws.Range("B7").Value = 866627
$ws.Range("C7").Value = 5204
$ws.Range("D7").Value = 669026
$ws.Range("E7").Value = 183111
$ws.Range("G7").Value = 139
$ws.Range("H7").Value = 14490

# Filipinas moves up to row 26 with fresh data (was Indonesia's row)
$ws.Range("A26").Value = "Filipinas"
$ws.Range("B26").Value = 115980
$ws.Range("C26").Value = 3462
$ws.Range("D26").Value = 66270
$ws.Range("E26").Value = 47587
$ws.Range("G26").Value = 9
$ws.Range("H26").Value = 2123

# Indonesia drops to row 27, keeping its prior totals
$ws.Range("A27").Value = "Indonesia"
$ws.Range("B27").Value = 115056
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 72050
$ws.Range("E27").Value = 37618
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 5388

# Israel (row 36)
$ws.Range("B36").Value = 76763
$ws.Range("C36").Value = 565
$ws.Range("D36").Value = 51331
$ws.Range("E36").Value = 24868
$ws.Range("G36").Value = 3
$ws.Range("H36").Value = 564

# Singapur (row 46)
$ws.Range("B46").Value = 54254
$ws.Range("C46").Value = 908
$ws.Range("E46").Value = 6773

# Polonia (row 49)
$ws.Range("D49").Value = 35321
$ws.Range("E49").Value = 11090

# Barein (row 52)
$ws.Range("E52").Value = 2645
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 152

# Armenia (row 54)
$ws.Range("E54").Value = 7966
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 770

# Uzbekistan (row 61)
$ws.Range("E61").Value = 9096
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 167

# Eslovaquia (row 125)
$ws.Range("B125").Value = 2417
$ws.Range("C125").Value = 49
$ws.Range("D125").Value = 1777
$ws.Range("E125").Value = 611

# Estonia moves up to row 128 with fresh data (was Ruanda's row)
$ws.Range("A128").Value = "Estonia"
$ws.Range("B128").Value = 2113
$ws.Range("C128").Value = 22
$ws.Range("D128").Value = 1948
$ws.Range("E128").Value = 102
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 63

# Ruanda drops to row 129, keeping its prior totals
$ws.Range("A129").Value = "Ruanda"
$ws.Range("B129").Value = 2099
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 1222
$ws.Range("E129").Value = 872
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 5
